$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.720800174688467
$ws.Range("D2").Value = 0.09932303673510723

$ws.Range("C3").Value = -0.4287187145169435
$ws.Range("D3").Value = 0.6722965672392491

$ws.Range("C4").Value = 0.2439908710706419
$ws.Range("D4").Value = 0.8094988809203987

$ws.Range("C5").Value = 0.9358898091791082
$ws.Range("D5").Value = 0.3594903728900205

$ws.Range("C6").Value = -1.583860312090974
$ws.Range("D6").Value = 0.1274961754878616

$ws.Range("C7").Value = -0.9682552645633468
$ws.Range("D7").Value = 0.3434418120315914

$ws.Range("C8").Value = -0.4718571527917298
$ws.Range("D8").Value = 0.6416751699281096

$ws.Range("C9").Value = 0.7755936212397754
$ws.Range("D9").Value = 0.4462467258789351

$ws.Range("C10").Value = 1.357420086661181
$ws.Range("D10").Value = 0.1884105589662197

$ws.Range("C11").Value = 0.6713615659656638
$ws.Range("D11").Value = 0.5089788169397762
